$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values to match the latest cryptos snapshot.
# Columns B (Coin) and C (Link) are plain text; D (Price) and E (Volume/1h)
# are also stored as text in the source data (note the "." thousands
# separators and padded "%" strings), so we force a text number format
# before writing numeric-looking values and then restore the default
# "Normal" style so no stray formatting is introduced.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.664.94'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.513.89'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.66%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.89%  '
$ws.Range('E6').Value = '  -2.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.583'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.53%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.539'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.55%  '
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.61'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.899.83'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.56'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.521.38'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.860'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.651.90'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.91'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.68%  '
$ws.Range('E20').Value = '  -1.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.61'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.47'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('E25').Value = '  -2.52%  '
$ws.Range('E26').Value = '  -3.02%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').Value = '  +13.24%  '
$ws.Range('E29').Value = '  +0.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.19'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '155.65'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.70'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.36'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.09'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.34%  '
$ws.Range('E36').Value = '  -2.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.61'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.81%  '
$ws.Range('E38').Value = '  -1.32%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '24.30'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.85%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.120'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.92%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  -1.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.05'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.63%  '
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.036.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '84.58'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.95'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '74.28'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.756.25'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.190'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.39%  '
